$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.993.64'
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").Value = '3.145.78'
$ws.Range("E3").Value = '  +2.30%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.43'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +2.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.60'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +3.80%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = '3.144.41'
$ws.Range("E8").Value = '  +2.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +4.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +6.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.12'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -1.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.504'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +7.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +12.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.49'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +6.79%  '

$ws.Range("D15").Value = '3.660.05'
$ws.Range("E15").Value = '  +2.31%  '

$ws.Range("D16").Value = '65.075.10'
$ws.Range("E16").Value = '  +1.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.18'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +6.51%  '

$ws.Range("D18").Value = '3.153.09'
$ws.Range("E18").Value = '  +2.55%  '

$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '509.61'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +6.72%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.89'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +7.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.730'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +8.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.49'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +13.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.85'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +4.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.53'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +5.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.92'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +4.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.77'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +9.44%  '

$ws.Range("E29").Value = '  +6.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.96'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +7.09%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +0.16%  '

$ws.Range("E32").Value = '  +3.94%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.65'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +6.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.03'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +8.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.58'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +6.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.62'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -0.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '472.28'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +3.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0423'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +4.42%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0857'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +4.06%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.04'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +0.86%  '

$ws.Range("D41").Value = '3.111.54'
$ws.Range("E41").Value = '  +5.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.61'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +4.63%  '

$ws.Range("E43").Value = '  +4.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.291'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +11.81%  '

$ws.Range("E45").Value = '  +13.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.16'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +4.94%  '

$ws.Range("D47").Value = '0.0₃0577'
$ws.Range("E47").Value = '  +12.17%  '

$ws.Range("E49").Value = '  +3.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +11.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.61'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -2.20%  '
